# UPDATE data for Washington, D.C.
$wb = $excel.ActiveWorkbook

# "2025" sheet: base year, literal value update
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("E2").Value = 3300.0000000000005
$ws2025.Range("H2").Value = 5.5

# "2030" sheet: p_ates_c_inv now derived from 2025 value with a 20% decline factor
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("E2").Formula = "='2025'!E2*(1-0.25*0.2)"
$ws2030.Range("H2").Value = 5.5

# "2035" sheet: 40% decline factor
$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("E2").Formula = "='2025'!E2*(1-0.25*0.4)"
$ws2035.Range("H2").Value = 5.5

# "2040" sheet: 60% decline factor
$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("E2").Formula = "='2025'!E2*(1-0.25*0.6)"
$ws2040.Range("H2").Value = 5.5

# "2045" sheet: 80% decline factor
$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("E2").Formula = "='2025'!E2*(1-0.25*0.8)"
$ws2045.Range("H2").Value = 5.5

# "2050" sheet: 100% decline factor
$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("E2").Formula = "='2025'!E2*(1-0.25*1)"
$ws2050.Range("H2").Value = 5.5
